# Auto-generated: reorders/updates rows 2-13 in Sheet1 per the
# commit's weekly fruit/vegetable price-data shuffle (Zapallo italiano).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: was D=44284, now D=44243 (values from former row 10)
$ws.Range("D2").Value = 44243
$ws.Range("J2").Value = 80
$ws.Range("L2").Value = 11000
$ws.Range("M2").Value = 10375
$ws.Range("O2").Value = 'Provincia de Quillota'
$ws.Range("P2").Value = 173

# Row 3: was D=44277, now D=44315 (values from former row 11)
$ws.Range("D3").Value = 44315

# Row 4: was D=44291, now D=44405 (values from former row 5)
$ws.Range("D4").Value = 44405
$ws.Range("J4").Value = 45
$ws.Range("N4").Value = '$/caja 50 unidades'
$ws.Range("O4").Value = 'Provincia de Quillota'
$ws.Range("P4").Value = 180
$ws.Range("Q4").Value = 50

# Row 5: was D=44405, now D=44333 (values from former row 8)
$ws.Range("D5").Value = 44333
$ws.Range("J5").Value = 25
$ws.Range("K5").Value = 10000
$ws.Range("L5").Value = 11000
$ws.Range("M5").Value = 10400
$ws.Range("N5").Value = '$/caja 60 unidades'
$ws.Range("O5").Value = 'Provincia de Limarí'
$ws.Range("P5").Value = 173
$ws.Range("Q5").Value = 60

# Row 6: was D=44186, now D=44585 (values from former row 9)
$ws.Range("D6").Value = 44585
$ws.Range("J6").Value = 30
$ws.Range("K6").Value = 11000
$ws.Range("L6").Value = 11000
$ws.Range("M6").Value = 11000
$ws.Range("P6").Value = 183

# Row 8: was D=44333, now D=44277 (values from former row 3)
$ws.Range("D8").Value = 44277
$ws.Range("L8").Value = 10000
$ws.Range("M8").Value = 10000
$ws.Range("P8").Value = 167

# Row 9: was D=44585, now D=44186 (values from former row 6)
$ws.Range("D9").Value = 44186
$ws.Range("J9").Value = 15
$ws.Range("K9").Value = 7000
$ws.Range("L9").Value = 7000
$ws.Range("M9").Value = 7000
$ws.Range("P9").Value = 117

# Row 10: was D=44243, now D=44200 (values from former row 13)
$ws.Range("D10").Value = 44200
$ws.Range("J10").Value = 10
$ws.Range("K10").Value = 9000
$ws.Range("L10").Value = 9000
$ws.Range("M10").Value = 9000
$ws.Range("O10").Value = 'Provincia de Limarí'
$ws.Range("P10").Value = 150

# Row 11: was D=44315, now D=44179 (values from former row 12)
$ws.Range("D11").Value = 44179
$ws.Range("J11").Value = 15
$ws.Range("K11").Value = 7000
$ws.Range("L11").Value = 7000
$ws.Range("M11").Value = 7000
$ws.Range("P11").Value = 117

# Row 12: was D=44179, now D=44284 (values from former row 2)
$ws.Range("D12").Value = 44284
$ws.Range("J12").Value = 35
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = 10000
$ws.Range("P12").Value = 167

# Row 13: was D=44200, now D=44291 (values from former row 4)
$ws.Range("D13").Value = 44291
$ws.Range("J13").Value = 20

